# Adds a "Sources" tab (data-source references) as the new last sheet of the
# workbook, and makes it the active/selected sheet - mirroring the authored
# commit "added data source tabs in each Excel file".

$wb = $excel.ActiveWorkbook

# Insert the new worksheet after the current last sheet ("Fossil Fuel Denmark")
# so it lands at the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sources"

# General Motors source
$ws.Range("B2").Value = "General Motors 2017 annual report:"
$ws.Range("B3").Value = "https://www.gm.com/content/dam/gm/en_us/english/Group4/InvestorsPDFDocuments/02-pdfs/10-K.pdf"

# Fossil Fuel Denmark sources (note: write B6 before B5 to match the shared
# string order of the authored workbook)
$ws.Range("B6").Value = "www.statistikbanken.dk/ENE3H"
$ws.Range("B5").Value = "Fuel type by industry in Denmark:"

# Busiest container ports source
$ws.Range("B8").Value = "Busiest container ports: "
$ws.Range("B9").Value = "https://en.wikipedia.org/wiki/List_of_busiest_container_ports"

# Busiest airports source
$ws.Range("B11").Value = "Busiest airports by passenger traffic"
$ws.Range("B12").Value = "https://en.wikipedia.org/wiki/List_of_busiest_airports_by_passenger_traffic"

# Make the new sheet the active tab and set its selection, matching the
# recorded view state in the target workbook.
$ws.Activate()
$ws.Range("E24").Select() | Out-Null
